$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 566, shifting existing rows 566-607 down to 567-608.
$ws.Rows("566:566").Insert()

# Populate the newly inserted row 566 with the new data point for 2026/01/06.
# Force column A to be stored as plain text (matching the other date cells),
# not auto-converted into a date serial value, then restore the default
# "Normal" style so no stray number-format style is left attached to the cell.
$ws.Range("A566").NumberFormat = "@"
$ws.Range("A566").Value = "2026/01/06"
$ws.Range("A566").Style = "Normal"

$ws.Range("B566").Value = "火"
$ws.Range("C566").Value = 9
$ws.Range("D566").Value = 17
